$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. M23: TODO -> UNDERWAY (copy style+value from H14, which already has that exact style/value pairing)
$ws.Range("H14").Copy()
$ws.Range("M23").PasteSpecial(-4122)
$ws.Range("M23").Value = $ws.Range("H14").Value2

# 2. L27's text changes in place (shared string index 90 keeps its slot, text updates)
$ws.Range("L27").Value = "Needs to look smoother"

# 3. M27: TODO -> SATISFACTORY (copy style+value from M9, which already has that exact style/value pairing)
$ws.Range("M9").Copy()
$ws.Range("M27").PasteSpecial(-4122)
$ws.Range("M27").Value = $ws.Range("M9").Value2

# 4. M28: TODO -> SATISFACTORY
$ws.Range("M9").Copy()
$ws.Range("M28").PasteSpecial(-4122)
$ws.Range("M28").Value = $ws.Range("M9").Value2

# 5. New row 30: K30/L30 plain text, M30 styled like a TODO (copy from M22)
$ws.Range("K30").Value = "When reaches destination"
$ws.Range("L30").Value = "Change start to avatar location"
$ws.Range("M22").Copy()
$ws.Range("M30").PasteSpecial(-4122)
$ws.Range("M30").Value = $ws.Range("M22").Value2

# 6. Selection moves to M23
$ws.Range("M23").Select() | Out-Null
